# Add a new "2022-Q4" sheet with fund holdings data, positioned right after
# the "总计" (summary) sheet, and insert a corresponding summary row into
# the "总计" sheet.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)            # "总计"
$q3 = $wb.Worksheets.Item(2)                 # "2022-Q3" (currently 2nd sheet)

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet right after "总计".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q4"

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows: A = running index (0-based), B/C = text, D/E/F/G = text
# (numeric-looking strings, matching the rest of the workbook), H = number.
$rows = @(
    @{A=0; B="165310"; C="建信沪深300指数增强（LOF）A"; D="4.07"; E="92.93"; F="2.32"; G="0.0944"; H=7},
    @{A=1; B="015387"; C="中欧沪深300指数增强A";        D="2.61"; E="91.39"; F="2.68"; G="0.0699"; H=4},
    @{A=2; B="004008"; C="中融鑫思路灵活配置混合A";      D="1.02"; E="37.86"; F="2.22"; G="0.0226"; H=3},
    @{A=3; B="015388"; C="中欧沪深300指数增强C";        D="0.68"; E="91.39"; F="2.68"; G="0.0182"; H=4},
    @{A=4; B="004009"; C="中融鑫思路灵活配置混合C";      D="0.70"; E="37.86"; F="2.22"; G="0.0155"; H=3},
    @{A=5; B="014049"; C="中银远见成长混合A";           D="0.77"; E="62.02"; F="1.69"; G="0.0130"; H=4},
    @{A=6; B="009208"; C="建信沪深300指数增强（LOF）C";  D="0.18"; E="92.93"; F="2.32"; G="0.0042"; H=7},
    @{A=7; B="014050"; C="中银远见成长混合C";           D="0.05"; E="62.02"; F="1.69"; G="0.0008"; H=4}
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row.A
    $newSheet.Cells.Item($r, 2).Value = "'" + $row.B
    $newSheet.Cells.Item($r, 3).Value = $row.C
    $newSheet.Cells.Item($r, 4).Value = "'" + $row.D
    $newSheet.Cells.Item($r, 5).Value = "'" + $row.E
    $newSheet.Cells.Item($r, 6).Value = "'" + $row.F
    $newSheet.Cells.Item($r, 7).Value = "'" + $row.G
    $newSheet.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}

# Formatting to match the other per-quarter sheets: bold + thin border on
# the header row and on the running-index column A, centered / top-aligned.
foreach ($rng in @($newSheet.Range("B1:H1"), $newSheet.Range("A2:A9"))) {
    $rng.Font.Bold = $true
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
}

# ---------------------------------------------------------------------
# 2. Insert the corresponding "2022-Q4" row into the "总计" sheet, right
#    below the header, pushing the existing quarters down by one row.
#    (Range.Value *read* returns an opaque placeholder in this host, so
#    shift cell-by-cell using Value2 for reads and Value for writes,
#    working bottom-up so a row is read before it gets overwritten.)
# ---------------------------------------------------------------------
for ($srcRow = 8; $srcRow -ge 2; $srcRow--) {
    $dstRow = $srcRow + 1
    $b = $summary.Cells.Item($srcRow, 2).Value2
    $c = $summary.Cells.Item($srcRow, 3).Value2
    $d = $summary.Cells.Item($srcRow, 4).Value2
    $summary.Cells.Item($dstRow, 2).Value = $b
    $summary.Cells.Item($dstRow, 3).Value = $c
    $summary.Cells.Item($dstRow, 4).Value = $d
}

# Column A is just a 0-based running index; extend it to the new last row.
$summary.Cells.Item(9, 1).Value = 7

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 8
$summary.Range("D2").Value = 0.24

Write-Output "done"
